# Auto-generated: update Lamia Profits market-price/profit values across all leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2790.2856
$ws.Range("J29").Value = 9004
$ws.Range("L29").Value = 27012
$ws.Range("N29").Value = -27574
$ws.Range("H38").Value = 786
$ws.Range("I38").Value = 232.5
$ws.Range("K38").Value = 697.5
$ws.Range("M38").Value = -325.5
$ws.Range("H135").Value = 1675.1578
$ws.Range("I135").Value = 1114.875
$ws.Range("K135").Value = 10033.875
$ws.Range("M135").Value = -7498.875
$ws.Range("H138").Value = 4054.5312
$ws.Range("I138").Value = 2431.611
$ws.Range("J138").Value = 6141.143
$ws.Range("K138").Value = 7294.833
$ws.Range("L138").Value = 18423.429
$ws.Range("M138").Value = -2154.833
$ws.Range("N138").Value = -28703.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7739.981
$ws.Range("I61").Value = 7233.2
$ws.Range("K61").Value = 7233.2
$ws.Range("M61").Value = -7021.2
$ws.Range("H74").Value = 9261753
$ws.Range("I74").Value = 15874154
$ws.Range("J74").Value = 4391.933
$ws.Range("K74").Value = 15874154
$ws.Range("L74").Value = 4391.933
$ws.Range("M74").Value = -15873280
$ws.Range("N74").Value = -6139.933
$ws.Range("H77").Value = 9261753
$ws.Range("I77").Value = 15874154
$ws.Range("J77").Value = 4391.933
$ws.Range("K77").Value = 79370770
$ws.Range("L77").Value = 21959.665
$ws.Range("M77").Value = -79366402
$ws.Range("N77").Value = -30695.665
$ws.Range("H97").Value = 945.5
$ws.Range("I97").Value = 987.8421
$ws.Range("J97").Value = 141
$ws.Range("K97").Value = 987.8421
$ws.Range("L97").Value = 141
$ws.Range("M97").Value = -491.8421
$ws.Range("N97").Value = -1133
$ws.Range("H136").Value = 7739.981
$ws.Range("I136").Value = 7233.2
$ws.Range("K136").Value = 21699.6
$ws.Range("M136").Value = -19149.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
$ws.Range("H98").Value = 72293
$ws.Range("J98").Value = 72293
$ws.Range("L98").Value = 72293
$ws.Range("N98").Value = -78283
$ws.Range("H105").Value = 24993.924
$ws.Range("I105").Value = 26491.445
$ws.Range("K105").Value = 26491.445
$ws.Range("M105").Value = -24744.445
$ws.Range("H133").Value = 77387.5
$ws.Range("J133").Value = 77387.5
$ws.Range("L133").Value = 77387.5
$ws.Range("N133").Value = -87507.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26330.805
$ws.Range("I31").Value = 1666.0526
$ws.Range("J31").Value = 43687.48
$ws.Range("K31").Value = 1666.0526
$ws.Range("L31").Value = 43687.48
$ws.Range("M31").Value = -1371.0526
$ws.Range("N31").Value = -44277.48
$ws.Range("H34").Value = 26330.805
$ws.Range("I34").Value = 1666.0526
$ws.Range("J34").Value = 43687.48
$ws.Range("K34").Value = 1666.0526
$ws.Range("L34").Value = 43687.48
$ws.Range("M34").Value = -1464.0526
$ws.Range("N34").Value = -44091.48
$ws.Range("H55").Value = 24000
$ws.Range("H58").Value = 3944.2
$ws.Range("I58").Value = 1703.3125
$ws.Range("J58").Value = 7928
$ws.Range("K58").Value = 1703.3125
$ws.Range("L58").Value = 7928
$ws.Range("M58").Value = -1500.3125
$ws.Range("N58").Value = -8334
$ws.Range("H122").Value = 4109.423
$ws.Range("I122").Value = 1676.5294
$ws.Range("K122").Value = 5029.5882
$ws.Range("M122").Value = -2579.5882
$ws.Range("H136").Value = 3944.2
$ws.Range("I136").Value = 1703.3125
$ws.Range("J136").Value = 7928
$ws.Range("K136").Value = 5109.9375
$ws.Range("L136").Value = 23784
$ws.Range("M136").Value = -2559.9375
$ws.Range("N136").Value = -28884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 59
$ws.Range("I13").Value = 58.75
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 176.25
$ws.Range("L13").Value = 180
$ws.Range("M13").Value = -8.25
$ws.Range("N13").Value = -516
$ws.Range("H50").Value = 66673570
$ws.Range("I50").Value = 166666720
$ws.Range("J50").Value = 11465.333
$ws.Range("K50").Value = 500000160
$ws.Range("L50").Value = 34395.999
$ws.Range("M50").Value = -499999679
$ws.Range("N50").Value = -35357.999
$ws.Range("H53").Value = 66673570
$ws.Range("I53").Value = 166666720
$ws.Range("J53").Value = 11465.333
$ws.Range("K53").Value = 500000160
$ws.Range("L53").Value = 34395.999
$ws.Range("M53").Value = -499999679
$ws.Range("N53").Value = -35357.999
$ws.Range("H54").Value = 2859.25
$ws.Range("I54").Value = 799.6667
$ws.Range("J54").Value = 4095
$ws.Range("K54").Value = 2399.0001
$ws.Range("L54").Value = 12285
$ws.Range("M54").Value = -1840.0001
$ws.Range("N54").Value = -13403

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 42497.5
$ws.Range("I34").Value = 39995
$ws.Range("K34").Value = 39995
$ws.Range("M34").Value = -39727
$ws.Range("H44").Value = 23330
$ws.Range("J44").Value = 24995
$ws.Range("L44").Value = 24995
$ws.Range("N44").Value = -26187
$ws.Range("H70").Value = 9992.583000000001
$ws.Range("I70").Value = 8406.814
$ws.Range("J70").Value = 14749.889
$ws.Range("K70").Value = 8406.814
$ws.Range("L70").Value = 14749.889
$ws.Range("M70").Value = -8136.814
$ws.Range("N70").Value = -15289.889
$ws.Range("H73").Value = 9992.583000000001
$ws.Range("I73").Value = 8406.814
$ws.Range("J73").Value = 14749.889
$ws.Range("K73").Value = 8406.814
$ws.Range("L73").Value = 14749.889
$ws.Range("M73").Value = -7470.814
$ws.Range("N73").Value = -16621.889
$ws.Range("H76").Value = 42497.5
$ws.Range("I76").Value = 39995
$ws.Range("K76").Value = 39995
$ws.Range("M76").Value = -39680
$ws.Range("H79").Value = 42497.5
$ws.Range("I79").Value = 39995
$ws.Range("K79").Value = 39995
$ws.Range("M79").Value = -38903
$ws.Range("H97").Value = 1477.1538
$ws.Range("J97").Value = 900
$ws.Range("L97").Value = 900
$ws.Range("N97").Value = -1892

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 662.2941
$ws.Range("I16").Value = 693.25
$ws.Range("J16").Value = 588
$ws.Range("K16").Value = 693.25
$ws.Range("L16").Value = 588
$ws.Range("M16").Value = -523.25
$ws.Range("H136").Value = 10097
$ws.Range("I136").Value = 5002
$ws.Range("K136").Value = 15006
$ws.Range("M136").Value = -12456

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11638.111
$ws.Range("I132").Value = 7058
$ws.Range("J132").Value = 13928.167
$ws.Range("K132").Value = 21174
$ws.Range("L132").Value = 41784.501
$ws.Range("M132").Value = -18644
$ws.Range("N132").Value = -46844.501
$ws.Range("H136").Value = 2890.5
$ws.Range("I136").Value = 2042.3684
$ws.Range("K136").Value = 6127.1052
$ws.Range("M136").Value = -3577.1052
